$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 51.59157666666666
$ws.Range("H2").Value = 154.77473
$ws.Range("I2").Value = 0.2641250550177587
$ws.Range("J2").Value = 0.2641250550177588
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 218.7785543333333
$ws.Range("N2").Value = 656.3356630000001
$ws.Range("O2").Value = 0.7837094150017259
$ws.Range("P2").Value = 0.7837094150017259
$ws.Range("Q2").Value = 11287.13055891067
$ws.Range("R2").Value = 101584.175030196
$ws.Range("S2").Value = 0.2069972923552663
$ws.Range("T2").Value = 0.2069972923552664

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 51.59157666666666
$ws.Range("H3").Value = 154.77473
$ws.Range("I3").Value = 0.2641250550177587
$ws.Range("J3").Value = 0.2641250550177588
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 46.29469433333333
$ws.Range("N3").Value = 138.884083
$ws.Range("O3").Value = 0.1658370397602197
$ws.Range("P3").Value = 0.1658370397602197
$ws.Range("Q3").Value = 2388.416271958065
$ws.Range("R3").Value = 21495.74644762259
$ws.Range("S3").Value = 0.04380171725065026
$ws.Range("T3").Value = 0.04380171725065027

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 51.59157666666666
$ws.Range("H4").Value = 154.77473
$ws.Range("I4").Value = 0.2641250550177587
$ws.Range("J4").Value = 0.2641250550177588
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.978882
$ws.Range("N4").Value = 20.936646
$ws.Range("O4").Value = 0.02499977909741928
$ws.Range("P4").Value = 0.02499977909741927
$ws.Range("Q4").Value = 360.0515257506199
$ws.Range("R4").Value = 3240.463731755579
$ws.Range("S4").Value = 0.006603068029537681
$ws.Range("T4").Value = 0.006603068029537682

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 51.59157666666666
$ws.Range("H5").Value = 154.77473
$ws.Range("I5").Value = 0.2641250550177587
$ws.Range("J5").Value = 0.2641250550177588
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 7.105616
$ws.Range("N5").Value = 21.316848
$ws.Range("O5").Value = 0.02545376614063513
$ws.Range("P5").Value = 0.02545376614063513
$ws.Range("Q5").Value = 366.5899326278933
$ws.Range("R5").Value = 3299.309393651039
$ws.Range("S5").Value = 0.006722977382304419
$ws.Range("T5").Value = 0.006722977382304419

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 19.32115333333334
$ws.Range("H6").Value = 57.96346000000001
$ws.Range("I6").Value = 0.09891538535728452
$ws.Range("J6").Value = 0.09891538535728453
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 218.7785543333333
$ws.Range("N6").Value = 656.3356630000001
$ws.Range("O6").Value = 0.7837094150017259
$ws.Range("P6").Value = 0.7837094150017259
$ws.Range("Q6").Value = 4227.053994319333
$ws.Range("R6").Value = 38043.48594887399
$ws.Range("S6").Value = 0.07752091879302773
$ws.Range("T6").Value = 0.07752091879302773

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 19.32115333333334
$ws.Range("H7").Value = 57.96346000000001
$ws.Range("I7").Value = 0.09891538535728452
$ws.Range("J7").Value = 0.09891538535728453
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 46.29469433333333
$ws.Range("N7").Value = 138.884083
$ws.Range("O7").Value = 0.1658370397602197
$ws.Range("P7").Value = 0.1658370397602197
$ws.Range("Q7").Value = 894.4668877341313
$ws.Range("R7").Value = 8050.201989607182
$ws.Range("S7").Value = 0.01640383469439344
$ws.Range("T7").Value = 0.01640383469439344

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 19.32115333333334
$ws.Range("H8").Value = 57.96346000000001
$ws.Range("I8").Value = 0.09891538535728452
$ws.Range("J8").Value = 0.09891538535728453
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 6.978882
$ws.Range("N8").Value = 20.936646
$ws.Range("O8").Value = 0.02499977909741928
$ws.Range("P8").Value = 0.02499977909741927
$ws.Range("Q8").Value = 134.84004921724
$ws.Range("R8").Value = 1213.56044295516
$ws.Range("S8").Value = 0.002472862783268214
$ws.Range("T8").Value = 0.002472862783268214

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 19.32115333333334
$ws.Range("H9").Value = 57.96346000000001
$ws.Range("I9").Value = 0.09891538535728452
$ws.Range("J9").Value = 0.09891538535728453
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.105616
$ws.Range("N9").Value = 21.316848
$ws.Range("O9").Value = 0.02545376614063513
$ws.Range("P9").Value = 0.02545376614063513
$ws.Range("Q9").Value = 137.2886962637867
$ws.Range("R9").Value = 1235.59826637408
$ws.Range("S9").Value = 0.002517769086595125
$ws.Range("T9").Value = 0.002517769086595125

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 112.3724673333333
$ws.Range("H10").Value = 337.117402
$ws.Range("I10").Value = 0.5752951554216499
$ws.Range("J10").Value = 0.57529515542165
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 218.7785543333333
$ws.Range("N10").Value = 656.3356630000001
$ws.Range("O10").Value = 0.7837094150017259
$ws.Range("P10").Value = 0.7837094150017259
$ws.Range("Q10").Value = 24584.68595005639
$ws.Range("R10").Value = 221262.1735505075
$ws.Range("S10").Value = 0.4508642297088282
$ws.Range("T10").Value = 0.4508642297088283

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 112.3724673333333
$ws.Range("H11").Value = 337.117402
$ws.Range("I11").Value = 0.5752951554216499
$ws.Range("J11").Value = 0.57529515542165
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 46.29469433333333
$ws.Range("N11").Value = 138.884083
$ws.Range("O11").Value = 0.1658370397602197
$ws.Range("P11").Value = 0.1658370397602197
$ws.Range("Q11").Value = 5202.249026679151
$ws.Range("R11").Value = 46820.24124011237
$ws.Range("S11").Value = 0.0954052455635219
$ws.Range("T11").Value = 0.09540524556352192

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 112.3724673333333
$ws.Range("H12").Value = 337.117402
$ws.Range("I12").Value = 0.5752951554216499
$ws.Range("J12").Value = 0.57529515542165
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 6.978882
$ws.Range("N12").Value = 20.936646
$ws.Range("O12").Value = 0.02499977909741928
$ws.Range("P12").Value = 0.02499977909741927
$ws.Range("Q12").Value = 784.2341895681878
$ws.Range("R12").Value = 7058.107706113691
$ws.Range("S12").Value = 0.01438225180135674
$ws.Range("T12").Value = 0.01438225180135674

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 112.3724673333333
$ws.Range("H13").Value = 337.117402
$ws.Range("I13").Value = 0.5752951554216499
$ws.Range("J13").Value = 0.57529515542165
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 7.105616
$ws.Range("N13").Value = 21.316848
$ws.Range("O13").Value = 0.02545376614063513
$ws.Range("P13").Value = 0.02545376614063513
$ws.Range("Q13").Value = 798.4756018432106
$ws.Range("R13").Value = 7186.280416588896
$ws.Range("S13").Value = 0.01464342834794302
$ws.Range("T13").Value = 0.01464342834794302

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 12.044915
$ws.Range("H14").Value = 36.134745
$ws.Range("I14").Value = 0.06166440420330686
$ws.Range("J14").Value = 0.06166440420330688
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 218.7785543333333
$ws.Range("N14").Value = 656.3356630000001
$ws.Range("O14").Value = 0.7837094150017259
$ws.Range("P14").Value = 0.7837094150017259
$ws.Range("Q14").Value = 2635.169090767882
$ws.Range("R14").Value = 23716.52181691094
$ws.Range("S14").Value = 0.04832697414460359
$ws.Range("T14").Value = 0.04832697414460359

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 12.044915
$ws.Range("H15").Value = 36.134745
$ws.Range("I15").Value = 0.06166440420330686
$ws.Range("J15").Value = 0.06166440420330688
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 46.29469433333333
$ws.Range("N15").Value = 138.884083
$ws.Range("O15").Value = 0.1658370397602197
$ws.Range("P15").Value = 0.1658370397602197
$ws.Range("Q15").Value = 557.6156581959817
$ws.Range("R15").Value = 5018.540923763835
$ws.Range("S15").Value = 0.01022624225165406
$ws.Range("T15").Value = 0.01022624225165406

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 12.044915
$ws.Range("H16").Value = 36.134745
$ws.Range("I16").Value = 0.06166440420330686
$ws.Range("J16").Value = 0.06166440420330688
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 6.978882
$ws.Range("N16").Value = 20.936646
$ws.Range("O16").Value = 0.02499977909741928
$ws.Range("P16").Value = 0.02499977909741927
$ws.Range("Q16").Value = 84.06004048503
$ws.Range("R16").Value = 756.5403643652701
$ws.Range("S16").Value = 0.001541596483256644
$ws.Range("T16").Value = 0.001541596483256644

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 12.044915
$ws.Range("H17").Value = 36.134745
$ws.Range("I17").Value = 0.06166440420330686
$ws.Range("J17").Value = 0.06166440420330688
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 7.105616
$ws.Range("N17").Value = 21.316848
$ws.Range("O17").Value = 0.02545376614063513
$ws.Range("P17").Value = 0.02545376614063513
$ws.Range("Q17").Value = 85.58654074264001
$ws.Range("R17").Value = 770.27886668376
$ws.Range("S17").Value = 0.001569591323792571
$ws.Range("T17").Value = 0.001569591323792571
